# Update cryptocurrency price (D) and 1h volume-change (E) columns
# to the latest scraped snapshot. Price cells are forced to Text so
# Excel does not reinterpret values like "216.90" or "1.320" as numbers
# (which would silently drop the trailing zero / change the cell type).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.092.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.669.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5122'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.73%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2641'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06435'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07429'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.672.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.516'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5822'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008587'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.176.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.937'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.223'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.646'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1204'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06401'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.304'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.320'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.529'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.516'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.638'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.020'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6104'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.361'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.650'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.198'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01611'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.081.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8626'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.816.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("E45").Value = '  +6.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.071'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05206'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4287'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.957'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.17%  '
